$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = "制作Model部分的类结构图，并与Kyle Review"
$ws.Range("D3").Value = "Open"

$ws.Range("D4").Select()
